# Update column F ("想去人数") values on several worksheets to match the
# latest scrape output (gh-pages data refresh at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> (row -> new value for column F)
$updates = @{
    "展览" = @{
        4 = 497
        5 = 2296
        6 = 3
        12 = 72
        16 = 681
        17 = 175
        19 = 7405
        20 = 8273
        36 = 236
        43 = 361
        47 = 197
        48 = 179
        49 = 21
    }
    "演出" = @{
        2 = 20
        3 = 37
        5 = 66
    }
    "本地生活" = @{
        4 = 293
        5 = 146
    }
    "全部类型" = @{
        5 = 20
        6 = 146
        7 = 37
        8 = 497
        9 = 2296
        14 = 72
        16 = 681
        18 = 175
        19 = 66
        21 = 7405
        22 = 8273
        30 = 236
        43 = 361
        47 = 197
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $updates[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $ws.Range("F$row").Value = $rowMap[$row]
    }
}
